$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '66.420.57'
$ws.Range('E2').Value = '  +4.09%  '

$ws.Range('D3').Value = '3.417.24'
$ws.Range('E3').Value = '  +3.83%  '

$ws.Range('D4').Value = '0.999'
$ws.Range('E4').Value = '  -0.24%  '

$ws.Range('D5').Value = '184.37'
$ws.Range('E5').Value = '  +3.59%  '

$ws.Range('D6').Value = '546.48'
$ws.Range('E6').Value = '  +4.13%  '

$ws.Range('D7').Value = '0.612'
$ws.Range('E7').Value = '  +1.66%  '

$ws.Range('D8').Value = '3.412.26'
$ws.Range('E8').Value = '  +3.78%  '

$ws.Range('E9').Value = '  -0.15%  '

$ws.Range('D10').Value = '0.635'
$ws.Range('E10').Value = '  +4.43%  '

$ws.Range('D11').Value = '56.04'
$ws.Range('E11').Value = '  -3.07%  '

$ws.Range('D12').Value = '0.149'
$ws.Range('E12').Value = '  +12.18%  '

$ws.Range('E13').Value = '  +5.86%  '

$ws.Range('D14').Value = '9.36'
$ws.Range('E14').Value = '  +3.15%  '

$ws.Range('D15').Value = '3.945.26'
$ws.Range('E15').Value = '  +3.03%  '

$ws.Range('E16').Value = '  +3.34%  '

$ws.Range('D17').Value = '3.406.46'
$ws.Range('E17').Value = '  +3.27%  '

$ws.Range('D18').Value = '18.21'
$ws.Range('E18').Value = '  +4.66%  '

$ws.Range('D19').Value = '66.582.87'
$ws.Range('E19').Value = '  +4.36%  '

$ws.Range('E20').Value = '  +5.22%  '

$ws.Range('E21').Value = '  +4.30%  '

$ws.Range('D22').Value = '404.90'
$ws.Range('E22').Value = '  +8.46%  '

$ws.Range('D23').Value = '11.98'
$ws.Range('E23').Value = '  +8.21%  '

$ws.Range('D24').Value = '4.26'
$ws.Range('E24').Value = '  +9.32%  '

$ws.Range('B25').Value = 'PancakeSwap'
$ws.Range('C25').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D25').Value = '3.87'
$ws.Range('E25').Value = '  +2.44%  '

$ws.Range('B26').Value = 'Litecoin'
$ws.Range('C26').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('D26').Value = '84.07'
$ws.Range('E26').Value = '  +4.52%  '

$ws.Range('B27').Value = 'LEO'
$ws.Range('C27').Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range('D27').Value = '6.22'
$ws.Range('E27').Value = '  +2.52%  '

$ws.Range('B28').Value = 'ImmutableX'
$ws.Range('C28').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D28').Value = '2.88'
$ws.Range('E28').Value = '  +7.79%  '

$ws.Range('D29').Value = '11.68'
$ws.Range('E29').Value = '  +2.85%  '

$ws.Range('E30').Value = '  +3.09%  '

$ws.Range('D31').Value = '30.02'
$ws.Range('E31').Value = '  +4.12%  '

$ws.Range('D32').Value = '668.28'
$ws.Range('E32').Value = '  +3.88%  '

$ws.Range('D33').Value = '6.82'
$ws.Range('E33').Value = '  +3.20%  '

$ws.Range('E34').Value = '  +2.99%  '

$ws.Range('E35').Value = '  +4.16%  '

$ws.Range('D36').Value = '58.81'
$ws.Range('E36').Value = '  -0.62%  '

$ws.Range('B37').Value = 'PEPE'
$ws.Range('C37').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D37').Value = '0.0₃0819'
$ws.Range('E37').Value = '  +17.31%  '

$ws.Range('B38').Value = 'InjectiveProtocol'
$ws.Range('C38').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D38').Value = '38.49'
$ws.Range('E38').Value = '  +5.24%  '

$ws.Range('D39').Value = '0.402'
$ws.Range('E39').Value = '  +3.35%  '

$ws.Range('D40').Value = '1.00'

$ws.Range('D41').Value = '2.80'
$ws.Range('E41').Value = '  +13.95%  '

$ws.Range('E42').Value = '  +6.60%  '

$ws.Range('E43').Value = '  +18.97%  '

$ws.Range('D44').Value = '0.998'
$ws.Range('E44').Value = '  -0.09%  '

$ws.Range('D45').Value = '3.031.09'
$ws.Range('E45').Value = '  +2.97%  '

$ws.Range('D46').Value = '2.93'
$ws.Range('E46').Value = '  +8.79%  '

$ws.Range('B47').Value = 'ApeXProtocol'
$ws.Range('C47').Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range('D47').Value = '3.28'
$ws.Range('E47').Value = '  +7.46%  '

$ws.Range('B48').Value = 'VeChain'
$ws.Range('C48').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D48').Value = '0.0417'
$ws.Range('E48').Value = '  +4.99%  '

$ws.Range('E49').Value = '  +3.83%  '

$ws.Range('E50').Value = '  +12.15%  '

$ws.Range('B51').Value = 'Stellar'
$ws.Range('C51').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D51').Value = '0.129'
$ws.Range('E51').Value = '  +2.99%  '
